# The paragraph originally ends with a single run whose text is:
#   ". All will extend the Airport class. "
# The edit splits that one run into three runs with identical
# run-formatting (rFonts Times New Roman + lang en-US), but with the
# run-level rsid attribute dropped (the runs were effectively retyped):
#   "."  " "  "All will extend the Airport class. "

$d = $word.ActiveDocument

$targetText = ". All will extend the Airport class. "

# Step 1: re-assert the same text through Find/Replace. This forces the
# engine to mint a fresh run for the match (dropping the legacy
# w:rsidRPr carried on the original run) without altering the visible
# text or formatting.
$rng = $d.Content
$replaced = $rng.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, $targetText, 2)
if (-not $replaced) {
    throw "Could not find target text to replace"
}

# Step 2: locate the (now rsid-less) run again so we know the exact
# character offsets to split on.
$rng2 = $d.Content
$found = $rng2.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if (-not $found) {
    throw "Could not find target text to split"
}
$start = $rng2.Start

# Step 3: force a run split at the boundaries between "." | " " | "All...".
# Toggling (and immediately reverting) a character formatting property on
# the single space character (offsets start+1 .. start+2) makes the
# engine break the run at both of those offsets, without changing the
# resulting formatting at all, yielding three runs:
#   [start, start+1)   -> "."
#   [start+1, start+2) -> " "
#   [start+2, end)      -> "All will extend the Airport class. "
$splitRange = $d.Range($start + 1, $start + 2)
$splitRange.Bold = 1
$splitRange = $d.Range($start + 1, $start + 2)
$splitRange.Bold = 0

Write-Output "Done splitting run."
